$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Add the Product module: two new worksheets, "PRODUCTS" and
# "SERVICES", inserted right after "org" and before "product" so the
# final tab order becomes:
#   contact, org, PRODUCTS, SERVICES, product, Sheet1
# ------------------------------------------------------------------
$orgSheet = $wb.Worksheets.Item("org")
$productsSheet = $wb.Worksheets.Add([System.Type]::Missing, $orgSheet)
$productsSheet.Name = "PRODUCTS"
$servicesSheet = $wb.Worksheets.Add([System.Type]::Missing, $productsSheet)
$servicesSheet.Name = "SERVICES"

# --- PRODUCTS sheet data ---
$productsSheet.Range("A1").Value = "TC_ID"
$productsSheet.Range("B1").Value = "TestcaseName"
$productsSheet.Range("C1").Value = "ProductRef"
$productsSheet.Range("D1").Value = "Label"
$productsSheet.Range("E1").Value = "Weight"
$productsSheet.Range("F1").Value = "Units"

$productsSheet.Range("A2").Value = "tc_01"
$productsSheet.Range("B2").Value = "modifyProductTest"
$productsSheet.Range("C2").Value = "LG_Oven"
$productsSheet.Range("D2").Value = "Oven"
$productsSheet.Range("E2").Value = "'5"
$productsSheet.Range("F2").Value = "mg"

$productsSheet.Range("A3").Value = "tc_02"
$productsSheet.Range("B3").Value = "deleteServiceTest"
$productsSheet.Range("C3").Value = "Car_Service"
$productsSheet.Range("D3").Value = "Car Wash"

# --- SERVICES sheet data ---
$servicesSheet.Range("A1").Value = "TC_ID"
$servicesSheet.Range("B1").Value = "TestcaseName"
$servicesSheet.Range("C1").Value = "ProductRef"
$servicesSheet.Range("D1").Value = "Label"

$servicesSheet.Range("A2").Value = "tc_02"
$servicesSheet.Range("B2").Value = "deleteServiceTest"
$servicesSheet.Range("C2").Value = "Car_Service"
$servicesSheet.Range("D2").Value = "Car Wash"

# ------------------------------------------------------------------
# Selections / active sheet, matching the saved workbook view state:
# "contact" keeps a block selection, "PRODUCTS" selects its whole
# used range, and "SERVICES" ends up as the active/selected tab.
# ------------------------------------------------------------------
$contactSheet = $wb.Worksheets.Item("contact")
$contactSheet.Activate()
$contactSheet.Range("A1:E2").Select() | Out-Null

$productsSheet.Activate()
$productsSheet.Range("A1:F3").Select() | Out-Null

$servicesSheet.Activate()
$servicesSheet.Range("E1:F1").Select() | Out-Null
